# #327 Ajout des profils d'acces a58d18c1e8091c98efec92c8c093b361a253eee5
#
# 1. Bump the "Date" metadata value on the Metadata sheet.
# 2. On the Elements sheet, swap the content of the two "Mapping" columns
#    (AK = col 37, AL = col 38) for the rows that actually differ, and
#    swap their column widths to match.

$wb = $excel.ActiveWorkbook
$wsMeta = $wb.Worksheets.Item("Metadata")
$wsElem = $wb.Worksheets.Item("Elements")

# 1. Metadata!B8 holds the "Date" property value.
$wsMeta.Range("B8").Value = "2024-03-19T13:17:15+00:00"

# 2. Swap the "Mapping: RIM Mapping" / "Mapping: Spécification métier vers
#    l'extension ROR CoordinateReliability" columns (AK <-> AL), row by row,
#    for the header and the rows whose mapping values differ.
$rows = @(1, 3, 5, 6)
foreach ($r in $rows) {
  $akCell = $wsElem.Cells.Item($r, 37)
  $alCell = $wsElem.Cells.Item($r, 38)
  $akVal = $akCell.Text
  $alVal = $alCell.Text
  $akCell.Value = $alVal
  $alCell.Value = $akVal
}

# 3. Swap the column widths that went along with the column re-order
#    (AK was 24.98046875 / AL was 76.828125 -> now the reverse).
#    ColumnWidth is stored internally with a +5/6 offset versus the raw
#    OOXML column width, so subtract that back out here.
$wsElem.Columns(37).ColumnWidth = 76.828125 - (5 / 6)
$wsElem.Columns(38).ColumnWidth = 24.98046875 - (5 / 6)
